$wb = $excel.ActiveWorkbook

# 1. Update status text "Ready for handoff" -> "In Translation" everywhere it occurs
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share this text)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# 2. Narrow the Status-related columns to match the new, shorter text.
#    Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
